# Fruta / hortaliza, semanal
# Insert a new weekly record at row 16 (pushing existing rows 16-28 down to 17-29)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16:28 down by one row to make room for the new record.
$ws.Rows.Item(16).Insert()

# Copy the date cell style (style "2" / date number format) from the row below
# so the new date cell renders correctly, then fill in the new row's values.
$ws.Range("D17").Copy()
$ws.Range("D16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44658
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100101
$ws.Cells.Item(16, 8).Value = "Berries"
$ws.Cells.Item(16, 9).Value = 100101001
$ws.Cells.Item(16, 10).Value = "Arándano (blue)"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 160
$ws.Cells.Item(16, 14).Value = 6500
$ws.Cells.Item(16, 15).Value = 7000
$ws.Cells.Item(16, 16).Value = 6750
$ws.Cells.Item(16, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Linares"
$ws.Cells.Item(16, 19).Value = 3375
$ws.Cells.Item(16, 20).Value = 2
